$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New daily TGA rows appended: rows 4352-4405 (dates 2023-01-16 .. 2023-03-27) ---
$newRows = @(
  @(44952, 572622, 568553, 446685, 568553),
  @(44953, 568553, 567827, 446685, 567827),
  @(44956, 567827, 579826, 446685, 579826),
  @(44957, 579826, 567908, 446685, 567908),
  @(44958, 567908, 500852, 567908, 500852),
  @(44959, 500852, 506179, 567908, 506179),
  @(44960, 506179, 477618, 567908, 477618),
  @(44963, 477618, 489948, 567908, 489948),
  @(44964, 489948, 528118, 567908, 528118),
  @(44965, 528118, 495838, 567908, 495838),
  @(44966, 495838, 494252, 567908, 494252),
  @(44967, 494252, 490201, 567908, 490201),
  @(44970, 490201, 501103, 567908, 501103),
  @(44971, 501103, 526992, 567908, 526992),
  @(44972, 526992, 439703, 567908, 439703),
  @(44973, 439703, 479039, 567908, 479039),
  @(44974, 479039, 475675, 567908, 475675),
  @(44978, 475675, 508286, 567908, 508286),
  @(44979, 508286, 451307, 567908, 451307),
  @(44980, 451307, 363666, 567908, 363666),
  @(44981, 363666, 381671, 567908, 381671),
  @(44984, 381671, 394017, 567908, 394017),
  @(44985, 394017, 415005, 567908, 415005),
  @(44986, 415005, 351015, 415005, 351015),
  @(44987, 351015, 355232, 415005, 355232),
  @(44988, 355232, 327193, 415005, 327193),
  @(44991, 327193, 340182, 415005, 340182),
  @(44992, 340182, 344724, 415005, 344724),
  @(44993, 344724, 311731, 415005, 311731),
  @(44994, 311731, 246969, 415005, 246969),
  @(44995, 246969, 208074, 415005, 208074),
  @(44998, 208074, 227312, 415005, 227312),
  @(44999, 227312, 253921, 415005, 253921),
  @(45000, 253921, 277643, 415005, 277643),
  @(45001, 277643, 285108, 415005, 285108),
  @(45002, 285108, 280148, 415005, 280148),
  @(45005, 280148, 267101, 415005, 267101),
  @(45006, 267101, 224604, 415005, 224604),
  @(45007, 224604, 199856, 415005, 199856),
  @(45008, 199856, 192910, 415005, 192910),
  @(45009, 192910, 187365, 415005, 187365),
  @(45012, 187365, 200926, 415005, 200926),
  @(45013, 200926, 166348, 415005, 166348),
  @(45014, 166348, 162758, 415005, 162758),
  @(45015, 162758, 194336, 415005, 194336),
  @(45016, 194336, 177692, 415005, 177692),
  @(45019, 177692, 173105, 177692, 173105),
  @(45020, 173105, 140347, 177692, 140347),
  @(45021, 140347, 140688, 177692, 140688),
  @(45022, 140688, 112965, 177692, 112965),
  @(45023, 112965, 110822, 177692, 110822),
  @(45026, 110822, 125004, 177692, 125004),
  @(45027, 125004, 107469, 177692, 107469),
  @(45028, 107469, 86554, 177692, 86554)
)

$r = 4352
foreach ($row in $newRows) {
    $ws.Cells.Item(4351, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "Treasury General Account (TGA)"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $r = $r + 1
}

# --- Style cleanup on existing F column cells (s="3" -> default; F4158 s="4" -> new index without applyBorder) ---
for ($row = 4159; $row -le 4351; $row++) {
    $ws.Cells.Item($row, 6).ClearFormats()
}

$fixCell = $ws.Cells.Item(4158, 6)
$fixCell.HorizontalAlignment = -4131
$fixCell.HorizontalAlignment = -4108

# --- Selection / scroll position ---
$ws.Range("I4379").Select()
